$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update personal details in row 2
$ws.Range("A2").Value = "Abhijeet"
$ws.Range("B2").Value = "Singh"

# Mobile number and date-of-birth must stay plain text, matching the
# original shared-string (text) cell type, not be auto-converted to a
# number/date. Force text format, assign, then restore the default
# "Normal" style so no extra style gets attached to the cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "+918698567733"
$ws.Range("D2").Style = "Normal"

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "1997-12-22"
$ws.Range("F2").Style = "Normal"

$ws.Range("G2").Value = "abhijeet1234"

# Update the active cell of the selection (sqref stays A1:G5, active cell moves to A5)
$ws.Range("A1:G5").Select()
$ws.Range("A5").Activate()
